$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 819, shifting existing rows 819:860 down to 820:861
$ws.Rows.Item(819).Insert()

# Populate the newly inserted row with the new data point.
# Force column A to be stored as plain text (not auto-converted to a date
# serial) by applying a text number format before assigning the value, then
# restore the cell's style to match the plain, unstyled data cells used
# throughout the rest of the table.
$ws.Range("A819").NumberFormat = "@"
$ws.Range("A819").Value = "2026/02/20"
$ws.Range("A819").Style = $ws.Range("A2").Style

$ws.Range("B819").Value = "金"
$ws.Range("C819").Value = 10
$ws.Range("D819").Value = 55
